$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column BF (58) holds the game date as a literal text label like
# "2008-06-20". The label previously stored the wrong date
# ("6-20-2007-08") because of how the NBA stats source displayed it.
# Setting the cell format to Text before writing keeps Excel from
# re-interpreting the "YYYY-MM-DD" string as a date serial value, and
# ClearFormats() afterwards restores the cell to its original
# (unstyled) appearance so only the value itself changes.
for ($r = 2; $r -le 31; $r++) {
    $cell = $ws.Range("BF$r")
    $cell.NumberFormat = "@"
    $cell.Value = "2008-06-20"
    $cell.ClearFormats()
}
